$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.166.26'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.57%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.862.21'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.26%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7182'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.80%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.82%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07736'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.13%  '

$ws.Range("E9").Value = '  +0.49%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.94'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.33%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08242'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.69%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.878.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.15%  '

$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7147'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.38%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.203'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.214.51'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.820'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.22%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '242.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007772'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.122.20'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.98%  '

$ws.Range("E22").Value = '  -0.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.954'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.39%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.19%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1587'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +9.65%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.895'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.31%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.14'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.493'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.54%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.302'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.351'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.02%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.078'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.90%  '

$ws.Range("E33").Value = '  -0.70%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.172'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7274'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.82%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.677'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.21%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01848'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.694'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.154.94'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8992'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.70%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.095'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.68%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.05'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.21%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.64'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.26%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.018.42'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5278'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.62%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.759'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.271'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.10%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.860'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.89%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9997'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.18%  '

Write-Host "Applied all changes"